$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 from numeric 1 to text "Yes" -> introduces shared string "Yes" (idx 6)
$ws.Range("D2").Value = "Yes"

# Update row 3 (was a placeholder row) to become the Admin - user management row
$ws.Range("B3").Value = "Validate Admin Module"          # introduces idx 7
$ws.Range("C3").Value = "Validate user management"        # introduces idx 8
$ws.Range("D3").Value = "Yes"                              # reuses idx 6

# Add row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Validate Admin Module"            # reuses idx 7
$ws.Range("C4").Value = "Validate Job Module"              # introduces idx 9

# Add row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Validate Admin Module"            # reuses idx 7
$ws.Range("C5").Value = "Validate Oragnization Module"     # introduces idx 10

# Now set the "No" flags after Oragnization Module string so ordering matches target
$ws.Range("D4").Value = "No"                               # introduces idx 11
$ws.Range("D5").Value = "No"                                # reuses idx 11

# Set the active selection to match target state
$ws.Range("B6").Select()
